# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Femacal de La Calera" / Zanahoria
# at sheet row 145 (pushing the existing rows 145-190 down to 146-191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 145..190 down to 146..191, creating a blank row 145.
$ws.Rows("145:145").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(145, 1).Value = 3
$ws.Cells.Item(145, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(145, 3).Value = "Coquimbo"
$ws.Cells.Item(145, 4).Value = 44210
$ws.Cells.Item(145, 5).Value = 5
$ws.Cells.Item(145, 6).Value = 100114013
$ws.Cells.Item(145, 7).Value = "Zanahoria"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 250
$ws.Cells.Item(145, 11).Value = 8000
$ws.Cells.Item(145, 12).Value = 8500
$ws.Cells.Item(145, 13).Value = 8260
$ws.Cells.Item(145, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(145, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(145, 16).Value = 413
$ws.Cells.Item(145, 17).Value = 20
$ws.Cells.Item(145, 18).Value = "Hortaliza"
